$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = -13.253
$ws.Range("A9").Value = -21.831
$ws.Range("C12").Value = -11.555
$ws.Range("C14").Value = -12.89
$ws.Range("A18").Value = -22.27
$ws.Range("A20").Value = -20.558
$ws.Range("C26").Value = -13.131
$ws.Range("A27").Value = -21.854
$ws.Range("C27").Value = -13.684
$ws.Range("C29").Value = -11.736
$ws.Range("A35").Value = -19.823
$ws.Range("C37").Value = -13.426
$ws.Range("C38").Value = -13.393
$ws.Range("C51").Value = -11.591
$ws.Range("C52").Value = -11.175
$ws.Range("C55").Value = -13.752
$ws.Range("A69").Value = -21.652
$ws.Range("C69").Value = -11.312
$ws.Range("C70").Value = -12.03
$ws.Range("A76").Value = -20.306
$ws.Range("A78").Value = -19.854
$ws.Range("C81").Value = -13.418
$ws.Range("A82").Value = -21.984
$ws.Range("A83").Value = -21.641
$ws.Range("C83").Value = -13.042
$ws.Range("A93").Value = -21.56
$ws.Range("C102").Value = -13.419
